$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3324.2104
$ws.Range("I62").Value = 5262.5
$ws.Range("J62").Value = 1914.5454
$ws.Range("K62").Value = 5262.5
$ws.Range("L62").Value = 1914.5454
$ws.Range("M62").Value = -4638.5
$ws.Range("N62").Value = -3162.5454
$ws.Range("H65").Value = 3324.2104
$ws.Range("I65").Value = 5262.5
$ws.Range("J65").Value = 1914.5454
$ws.Range("K65").Value = 26312.5
$ws.Range("L65").Value = 9572.726999999999
$ws.Range("M65").Value = -23192.5
$ws.Range("N65").Value = -15812.727
$ws.Range("H96").Value = 55562252
$ws.Range("I96").Value = 4540.25
$ws.Range("J96").Value = 100008424
$ws.Range("K96").Value = 13620.75
$ws.Range("L96").Value = 300025272
$ws.Range("M96").Value = -12247.75
$ws.Range("N96").Value = -300028018
$ws.Range("H98").Value = 34019.668
$ws.Range("I98").Value = 1127
$ws.Range("J98").Value = 88840.78
$ws.Range("K98").Value = 1127
$ws.Range("L98").Value = 88840.78
$ws.Range("M98").Value = 371
$ws.Range("N98").Value = -91836.78
$ws.Range("H120").Value = 45759
$ws.Range("J120").Value = 45759
$ws.Range("L120").Value = 45759
$ws.Range("N120").Value = -55435
$ws.Range("H122").Value = 34019.668
$ws.Range("I122").Value = 1127
$ws.Range("J122").Value = 88840.78
$ws.Range("K122").Value = 3381
$ws.Range("L122").Value = 266522.34
$ws.Range("M122").Value = -931
$ws.Range("N122").Value = -271422.34
$ws.Range("H123").Value = 27333.334
$ws.Range("J123").Value = 27333.334
$ws.Range("L123").Value = 27333.334
$ws.Range("N123").Value = -37133.334
$ws.Range("H126").Value = 37930.668
$ws.Range("J126").Value = 37930.668
$ws.Range("L126").Value = 37930.668
$ws.Range("N126").Value = -47810.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2145.25
$ws.Range("I45").Value = 1938.375
$ws.Range("J45").Value = 2559
$ws.Range("K45").Value = 1938.375
$ws.Range("L45").Value = 2559
$ws.Range("M45").Value = -1561.375
$ws.Range("N45").Value = -3313
$ws.Range("H113").Value = 36425.43
$ws.Range("J113").Value = 36425.43
$ws.Range("L113").Value = 36425.43
$ws.Range("N113").Value = -45103.43
$ws.Range("H122").Value = 1950.6111
$ws.Range("I122").Value = 1983.7333
$ws.Range("K122").Value = 5951.199900000001
$ws.Range("M122").Value = -3501.199900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2527.2856
$ws.Range("J105").Value = 2585.7778
$ws.Range("L105").Value = 2585.7778
$ws.Range("N105").Value = -6079.7778
$ws.Range("H107").Value = 2357.2173
$ws.Range("I107").Value = 2276.647
$ws.Range("J107").Value = 2585.5
$ws.Range("K107").Value = 2276.647
$ws.Range("L107").Value = 2585.5
$ws.Range("M107").Value = -356.6469999999999
$ws.Range("N107").Value = -6425.5
$ws.Range("H108").Value = 47680
$ws.Range("J108").Value = 47680
$ws.Range("L108").Value = 47680
$ws.Range("N108").Value = -55360
$ws.Range("H134").Value = 4140.7095
$ws.Range("I134").Value = 1944.6522
$ws.Range("K134").Value = 5833.9566
$ws.Range("M134").Value = -3298.9566
$ws.Range("H139").Value = 35166.668
$ws.Range("J139").Value = 35166.668
$ws.Range("L139").Value = 35166.668
$ws.Range("N139").Value = -45446.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 46776
$ws.Range("J100").Value = 46776
$ws.Range("L100").Value = 46776
$ws.Range("N100").Value = -48940
$ws.Range("H111").Value = 41468
$ws.Range("J111").Value = 41468
$ws.Range("L111").Value = 41468
$ws.Range("N111").Value = -49648
$ws.Range("H112").Value = 35696.715
$ws.Range("J112").Value = 35696.715
$ws.Range("L112").Value = 35696.715
$ws.Range("N112").Value = -38650.715
$ws.Range("H116").Value = 48431.5
$ws.Range("J116").Value = 48431.5
$ws.Range("L116").Value = 48431.5
$ws.Range("N116").Value = -57609.5
$ws.Range("H133").Value = 18857.846
$ws.Range("J133").Value = 18857.846
$ws.Range("L133").Value = 18857.846
$ws.Range("N133").Value = -23917.846
$ws.Range("H137").Value = 42849.918
$ws.Range("J137").Value = 42849.918
$ws.Range("L137").Value = 42849.918
$ws.Range("N137").Value = -53049.918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3098.0981
$ws.Range("I5").Value = 8871.5
$ws.Range("J5").Value = 1321.6666
$ws.Range("K5").Value = 26614.5
$ws.Range("L5").Value = 3964.9998
$ws.Range("M5").Value = -26502.5
$ws.Range("N5").Value = -4188.9998
$ws.Range("H87").Value = 7192.3335
$ws.Range("I87").Value = 7192.3335
$ws.Range("K87").Value = 21577.0005
$ws.Range("M87").Value = -20329.0005
$ws.Range("H90").Value = 7192.3335
$ws.Range("I90").Value = 7192.3335
$ws.Range("K90").Value = 64731.0015
$ws.Range("M90").Value = -58491.0015
$ws.Range("H135").Value = 3098.0981
$ws.Range("I135").Value = 8871.5
$ws.Range("J135").Value = 1321.6666
$ws.Range("K135").Value = 79843.5
$ws.Range("L135").Value = 11894.9994
$ws.Range("M135").Value = -77308.5
$ws.Range("N135").Value = -16964.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1814.4814
$ws.Range("I113").Value = 1755.6875
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1755.6875
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 414.3125
$ws.Range("N113").Value = -6240
$ws.Range("H126").Value = 2356.353
$ws.Range("I126").Value = 2628.6667
$ws.Range("J126").Value = 2050
$ws.Range("K126").Value = 7886.000100000001
$ws.Range("L126").Value = 6150
$ws.Range("M126").Value = -5416.000100000001
$ws.Range("N126").Value = -11090
$ws.Range("H133").Value = 35619
$ws.Range("J133").Value = 35619
$ws.Range("L133").Value = 35619
$ws.Range("N133").Value = -45739
$ws.Range("H138").Value = 43648
$ws.Range("J138").Value = 43648
$ws.Range("L138").Value = 43648
$ws.Range("N138").Value = -53928

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 34015
$ws.Range("J133").Value = 34015
$ws.Range("L133").Value = 34015
$ws.Range("N133").Value = -39075
$ws.Range("H137").Value = 41736.57
$ws.Range("J137").Value = 41736.57
$ws.Range("L137").Value = 41736.57
$ws.Range("N137").Value = -51936.57
$ws.Range("H138").Value = 41327.8
$ws.Range("J138").Value = 41327.8
$ws.Range("L138").Value = 41327.8
$ws.Range("N138").Value = -51607.8
